$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: "apache" -> "APACHE" for rows 2 through 13
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = "APACHE"
}

# Column D: numeric Severity values -> text representation of the same value
$severity = @{
    3  = "9.8"
    4  = "7.5"
    5  = "7.5"
    6  = "9.8"
    7  = "7.5"
    8  = "8.8"
    9  = "8.8"
    10 = "9.8"
    11 = "7.5"
    12 = "4.3"
    13 = "9.8"
}

foreach ($r in $severity.Keys) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = "'" + $severity[$r]
    $cell.Style = "Normal"
}
